$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(6, 6).Value = 680
$ws1.Cells.Item(7, 6).Value = 1247
$ws1.Cells.Item(9, 6).Value = 842
$ws1.Cells.Item(10, 6).Value = 709
$ws1.Cells.Item(13, 6).Value = 375
$ws1.Cells.Item(15, 6).Value = 985
$ws1.Cells.Item(16, 6).Value = 11152
$ws1.Cells.Item(17, 6).Value = 639
$ws1.Cells.Item(22, 6).Value = 279
$ws1.Cells.Item(23, 6).Value = 1783
$ws1.Cells.Item(26, 6).Value = 493
$ws1.Cells.Item(29, 6).Value = 284
$ws1.Cells.Item(30, 6).Value = 197
$ws1.Cells.Item(31, 6).Value = 265
$ws1.Cells.Item(34, 6).Value = 20
$ws1.Cells.Item(37, 6).Value = 187

# Sheet: 演出 (Show)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(7, 6).Value = 140
$ws2.Cells.Item(8, 6).Value = 188
$ws2.Cells.Item(10, 6).Value = 244
$ws2.Cells.Item(14, 6).Value = 6
$ws2.Cells.Item(16, 6).Value = 317
$ws2.Cells.Item(21, 6).Value = 5

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(9, 6).Value = 680
$ws4.Cells.Item(10, 6).Value = 1247
$ws4.Cells.Item(13, 6).Value = 140
$ws4.Cells.Item(14, 6).Value = 842
$ws4.Cells.Item(15, 6).Value = 709
$ws4.Cells.Item(18, 6).Value = 985
$ws4.Cells.Item(19, 6).Value = 11152
$ws4.Cells.Item(20, 6).Value = 244
$ws4.Cells.Item(21, 6).Value = 639
$ws4.Cells.Item(24, 6).Value = 279
$ws4.Cells.Item(25, 6).Value = 1783
$ws4.Cells.Item(26, 6).Value = 493
$ws4.Cells.Item(31, 6).Value = 6
$ws4.Cells.Item(33, 6).Value = 317
$ws4.Cells.Item(34, 6).Value = 284
$ws4.Cells.Item(36, 6).Value = 197
$ws4.Cells.Item(37, 6).Value = 265
$ws4.Cells.Item(40, 6).Value = 20
$ws4.Cells.Item(46, 6).Value = 187
$ws4.Cells.Item(47, 6).Value = 5
